$wb = $excel.ActiveWorkbook

# Rename sheets
$wsPayment = $wb.Worksheets.Item("Payment")
$wsPayment.Name = "payments"

$wsSubmission = $wb.Worksheets.Item("Submission")
$wsSubmission.Name = "submissions"

# Select cells / ranges before switching active sheet
$wsSubmission.Range("V12").Select()

$wsPayment.Activate()
$wsPayment.Range("F16").Select()
